# Auto-generated edit script applying Yojimbo_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 347144.03
$ws.Range("I17").Value = 750
$ws.Range("J17").Value = 376010.22
$ws.Range("K17").Value = 2250
$ws.Range("L17").Value = 1128030.66
$ws.Range("M17").Value = -2082
$ws.Range("N17").Value = -1128366.66
$ws.Range("H76").Value = 2835.5557
$ws.Range("I76").Value = 2865
$ws.Range("J76").Value = 2600
$ws.Range("K76").Value = 2865
$ws.Range("L76").Value = 2600
$ws.Range("M76").Value = -2550
$ws.Range("N76").Value = -3230
$ws.Range("H79").Value = 2835.5557
$ws.Range("I79").Value = 2865
$ws.Range("J79").Value = 2600
$ws.Range("K79").Value = 2865
$ws.Range("L79").Value = 2600
$ws.Range("M79").Value = -1773
$ws.Range("N79").Value = -4784
$ws.Range("H98").Value = 1234.9375
$ws.Range("I98").Value = 986.5185
$ws.Range("J98").Value = 2576.4
$ws.Range("K98").Value = 986.5185
$ws.Range("L98").Value = 2576.4
$ws.Range("M98").Value = 511.4815
$ws.Range("N98").Value = -5572.4
$ws.Range("H106").Value = 2940.0435
$ws.Range("I106").Value = 2720.9092
$ws.Range("J106").Value = 3008.9143
$ws.Range("K106").Value = 2720.9092
$ws.Range("L106").Value = 3008.9143
$ws.Range("M106").Value = -2089.9092
$ws.Range("N106").Value = -4270.9143
$ws.Range("H107").Value = 867.25
$ws.Range("I107").Value = 696
$ws.Range("J107").Value = 1266.8334
$ws.Range("K107").Value = 696
$ws.Range("L107").Value = 1266.8334
$ws.Range("M107").Value = 1224
$ws.Range("N107").Value = -5106.8334
$ws.Range("H122").Value = 1234.9375
$ws.Range("I122").Value = 986.5185
$ws.Range("J122").Value = 2576.4
$ws.Range("K122").Value = 2959.5555
$ws.Range("L122").Value = 7729.200000000001
$ws.Range("M122").Value = -509.5554999999999
$ws.Range("N122").Value = -12629.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 894.65717
$ws.Range("I2").Value = 781.6
$ws.Range("J2").Value = 1177.3
$ws.Range("K2").Value = 781.6
$ws.Range("L2").Value = 1177.3
$ws.Range("M2").Value = -668.6
$ws.Range("N2").Value = -1403.3
$ws.Range("H32").Value = 3514.15
$ws.Range("I32").Value = 2867.5264
$ws.Range("J32").Value = 15800
$ws.Range("K32").Value = 2867.5264
$ws.Range("L32").Value = 15800
$ws.Range("M32").Value = -2580.5264
$ws.Range("N32").Value = -16374
$ws.Range("H61").Value = 2532.4412
$ws.Range("I61").Value = 1517.25
$ws.Range("K61").Value = 1517.25
$ws.Range("M61").Value = -1305.25
$ws.Range("H63").Value = 2456.0454
$ws.Range("I63").Value = 2237.2354
$ws.Range("J63").Value = 3200
$ws.Range("K63").Value = 2237.2354
$ws.Range("L63").Value = 3200
$ws.Range("M63").Value = -1551.2354
$ws.Range("N63").Value = -4572
$ws.Range("H66").Value = 2456.0454
$ws.Range("I66").Value = 2237.2354
$ws.Range("J66").Value = 3200
$ws.Range("K66").Value = 11186.177
$ws.Range("L66").Value = 16000
$ws.Range("M66").Value = -7754.177
$ws.Range("N66").Value = -22864
$ws.Range("H116").Value = 894.65717
$ws.Range("I116").Value = 781.6
$ws.Range("J116").Value = 1177.3
$ws.Range("K116").Value = 781.6
$ws.Range("L116").Value = 1177.3
$ws.Range("M116").Value = 1512.4
$ws.Range("N116").Value = -5765.3
$ws.Range("H122").Value = 1647.619
$ws.Range("I122").Value = 1512.5
$ws.Range("J122").Value = 1730.7693
$ws.Range("K122").Value = 4537.5
$ws.Range("L122").Value = 5192.3079
$ws.Range("M122").Value = -2087.5
$ws.Range("N122").Value = -10092.3079
$ws.Range("H136").Value = 2532.4412
$ws.Range("I136").Value = 1517.25
$ws.Range("K136").Value = 4551.75
$ws.Range("M136").Value = -2001.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 894.65717
$ws.Range("I3").Value = 781.6
$ws.Range("J3").Value = 1177.3
$ws.Range("K3").Value = 781.6
$ws.Range("L3").Value = 1177.3
$ws.Range("M3").Value = -667.6
$ws.Range("N3").Value = -1405.3
$ws.Range("H105").Value = 1579.0476
$ws.Range("I105").Value = 1623.2778
$ws.Range("J105").Value = 1313.6666
$ws.Range("K105").Value = 1623.2778
$ws.Range("L105").Value = 1313.6666
$ws.Range("M105").Value = 123.7221999999999
$ws.Range("N105").Value = -4807.6666
$ws.Range("H107").Value = 1257.5714
$ws.Range("I107").Value = 995.1
$ws.Range("J107").Value = 1913.75
$ws.Range("K107").Value = 995.1
$ws.Range("L107").Value = 1913.75
$ws.Range("M107").Value = 924.9
$ws.Range("N107").Value = -5753.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 828.7368
$ws.Range("I16").Value = 906.1
$ws.Range("J16").Value = 742.7778
$ws.Range("K16").Value = 906.1
$ws.Range("L16").Value = 742.7778
$ws.Range("M16").Value = -619.1
$ws.Range("N16").Value = -1316.7778
$ws.Range("H62").Value = 2422.2222
$ws.Range("K62").Value = 2400
$ws.Range("M62").Value = -1776
$ws.Range("H65").Value = 2422.2222
$ws.Range("K65").Value = 12000
$ws.Range("M65").Value = -8880
$ws.Range("H105").Value = 800
$ws.Range("I105").Value = 775
$ws.Range("J105").Value = 850
$ws.Range("K105").Value = 775
$ws.Range("L105").Value = 850
$ws.Range("M105").Value = 972
$ws.Range("N105").Value = -4344
$ws.Range("H107").Value = 656.04
$ws.Range("I107").Value = 685.8182
$ws.Range("J107").Value = 437.66666
$ws.Range("K107").Value = 685.8182
$ws.Range("L107").Value = 437.66666
$ws.Range("M107").Value = 1234.1818
$ws.Range("N107").Value = -4277.66666
$ws.Range("H113").Value = 828.7368
$ws.Range("I113").Value = 906.1
$ws.Range("J113").Value = 742.7778
$ws.Range("K113").Value = 906.1
$ws.Range("L113").Value = 742.7778
$ws.Range("M113").Value = 1263.9
$ws.Range("N113").Value = -5082.7778

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 927.087
$ws.Range("J131").Value = 1052.7894
$ws.Range("L131").Value = 3158.3682
$ws.Range("N131").Value = -13238.3682

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4670
$ws.Range("I70").Value = 4337.5
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 4337.5
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -4067.5
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 4670
$ws.Range("I73").Value = 4337.5
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 4337.5
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -3401.5
$ws.Range("N73").Value = -7872
$ws.Range("H80").Value = 2450
$ws.Range("J80").Value = 2230
$ws.Range("L80").Value = 2230
$ws.Range("N80").Value = -4226
$ws.Range("H83").Value = 2450
$ws.Range("J83").Value = 2230
$ws.Range("L83").Value = 11150
$ws.Range("N83").Value = -21134
$ws.Range("H102").Value = 1078.16
$ws.Range("I102").Value = 896.1905
$ws.Range("J102").Value = 2033.5
$ws.Range("K102").Value = 896.1905
$ws.Range("L102").Value = 2033.5
$ws.Range("M102").Value = 725.8095
$ws.Range("N102").Value = -5277.5
$ws.Range("H122").Value = 2335.2856
$ws.Range("I122").Value = 1112
$ws.Range("J122").Value = 2946.9285
$ws.Range("K122").Value = 3336
$ws.Range("L122").Value = 8840.7855
$ws.Range("M122").Value = -886
$ws.Range("N122").Value = -13740.7855
$ws.Range("H132").Value = 2403
$ws.Range("I132").Value = 2241.55
$ws.Range("J132").Value = 2899.7693
$ws.Range("K132").Value = 6724.650000000001
$ws.Range("L132").Value = 8699.3079
$ws.Range("M132").Value = -4194.650000000001
$ws.Range("N132").Value = -13759.3079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2670.3845
$ws.Range("I40").Value = 2496.6667
$ws.Range("J40").Value = 2819.2856
$ws.Range("K40").Value = 2496.6667
$ws.Range("L40").Value = 2819.2856
$ws.Range("M40").Value = -2360.6667
$ws.Range("N40").Value = -3091.2856
$ws.Range("H61").Value = 2837.5
$ws.Range("I61").Value = 800
$ws.Range("J61").Value = 3128.5715
$ws.Range("K61").Value = 800
$ws.Range("L61").Value = 3128.5715
$ws.Range("M61").Value = -598
$ws.Range("N61").Value = -3532.5715
$ws.Range("H68").Value = 2553.7708
$ws.Range("I68").Value = 1020.5
$ws.Range("J68").Value = 2693.1592
$ws.Range("K68").Value = 1020.5
$ws.Range("L68").Value = 2693.1592
$ws.Range("M68").Value = -271.5
$ws.Range("N68").Value = -4191.1592
$ws.Range("H71").Value = 2553.7708
$ws.Range("I71").Value = 1020.5
$ws.Range("J71").Value = 2693.1592
$ws.Range("K71").Value = 5102.5
$ws.Range("L71").Value = 13465.796
$ws.Range("M71").Value = -1358.5
$ws.Range("N71").Value = -20953.796
$ws.Range("H113").Value = 2837.5
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 3128.5715
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 3128.5715
$ws.Range("M113").Value = 1370
$ws.Range("N113").Value = -7468.5715
$ws.Range("H132").Value = 6144.341
$ws.Range("I132").Value = 3851.7144
$ws.Range("J132").Value = 10156.4375
$ws.Range("K132").Value = 11555.1432
$ws.Range("L132").Value = 30469.3125
$ws.Range("M132").Value = -9025.143199999999
$ws.Range("N132").Value = -35529.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 12240
$ws.Range("J11").Value = 12240
$ws.Range("L11").Value = 12240
$ws.Range("N11").Value = -12524
$ws.Range("H107").Value = 12757.125
$ws.Range("I107").Value = 14436.714
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 43310.142
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -41390.142
$ws.Range("N107").Value = -6840
$ws.Range("H113").Value = 264.44446
$ws.Range("I113").Value = 214
$ws.Range("J113").Value = 327.5
$ws.Range("K113").Value = 642
$ws.Range("L113").Value = 982.5
$ws.Range("M113").Value = 1528
$ws.Range("N113").Value = -5322.5
$ws.Range("H122").Value = 271292.9
$ws.Range("I122").Value = 400851.7
$ws.Range("J122").Value = 1378.75
$ws.Range("K122").Value = 1202555.1
$ws.Range("L122").Value = 4136.25
$ws.Range("M122").Value = -1200105.1
$ws.Range("N122").Value = -9036.25
$ws.Range("H132").Value = 1029.575
$ws.Range("I132").Value = 649.0476
$ws.Range("J132").Value = 1450.1578
$ws.Range("K132").Value = 1947.1428
$ws.Range("L132").Value = 4350.4734
$ws.Range("M132").Value = 582.8571999999999
$ws.Range("N132").Value = -9410.473399999999
$ws.Range("H136").Value = 745
$ws.Range("I136").Value = 620.58826
$ws.Range("J136").Value = 1273.75
$ws.Range("K136").Value = 1861.76478
$ws.Range("L136").Value = 3821.25
$ws.Range("M136").Value = 688.23522
$ws.Range("N136").Value = -8921.25

